$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the formula columns for every student row (2-11) ---------
# "Yazılı Ortalama" (average of the two written exams)
$ws.Range("E2:E11").Style = "Normal"
$ws.Range("E2").Formula = "=(B2+C2)/2"
$ws.Range("E3:E11").Formula = "=(B3+C3)/2"

# "Genel Ortalama" (70% written average + 30% practical)
$ws.Range("F2").Formula = "=(E2*70%+D2*30%)"
$ws.Range("F3:F10").Formula = "=(E3*70%+D3*30%)"
$ws.Range("F11").Formula = "=(E11*70%+D11*30%)"

# "Durum" (pass/fail based on the general average)
$ws.Range("G2").Formula = '=IF(F2<50,"kaldı","geçti")'
$ws.Range("G3:G11").Formula = '=IF(F3<50,"kaldı","geçti")'

# --- Student identity card on the right (Numara / Ad Soyad / Bölüm) ---
$ws.Range("L4").Value = 20215070019
$ws.Range("L5").Value = "Kübra Çabuk"
$ws.Range("L6").Value = "YBS"

# --- Turn on AutoFilter for the data table -----------------------------
$null = $ws.Range("A1:G11").AutoFilter()
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sayfa1!`$A`$1:`$G`$11")
$n.Visible = $false

# --- Update the active selection ---------------------------------------
$null = $ws.Range("G2").Select()
